# Replace the "working set of sequences" rows (2-33) on the active sheet.
# Each row holds: B = a random/sequence number, C = image path, D = German
# word, E = category (face/flower). The commit reshuffles which
# image/word/category triple goes with each row and rerolls the B number,
# while row 1 (headers) and column A (the 0-based index) stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row=2; B=103; C="face/face013.jpg"; D="betteln"; E="face" },
    @{ Row=3; B=18; C="face/face002.jpg"; D="tollen"; E="face" },
    @{ Row=4; B=124; C="flower/flower030.jpg"; D="dienen"; E="flower" },
    @{ Row=5; B=22; C="flower/flower008.jpg"; D="ändern"; E="flower" },
    @{ Row=6; B=93; C="face/face029.jpg"; D="prüfen"; E="face" },
    @{ Row=7; B=28; C="face/face004.jpg"; D="passen"; E="face" },
    @{ Row=8; B=69; C="flower/flower031.jpg"; D="heißen"; E="flower" },
    @{ Row=9; B=113; C="face/face024.jpg"; D="hassen"; E="face" },
    @{ Row=10; B=100; C="face/face025.jpg"; D="proben"; E="face" },
    @{ Row=11; B=58; C="face/face001.jpg"; D="herrschen"; E="face" },
    @{ Row=12; B=2; C="flower/flower003.jpg"; D="lassen"; E="flower" },
    @{ Row=13; B=21; C="flower/flower000.jpg"; D="atmen"; E="flower" },
    @{ Row=14; B=86; C="face/face018.jpg"; D="parken"; E="face" },
    @{ Row=15; B=109; C="face/face005.jpg"; D="meinen"; E="face" },
    @{ Row=16; B=71; C="flower/flower020.jpg"; D="küssen"; E="flower" },
    @{ Row=17; B=35; C="flower/flower017.jpg"; D="zögern"; E="flower" },
    @{ Row=18; B=29; C="flower/flower011.jpg"; D="spüren"; E="flower" },
    @{ Row=19; B=111; C="face/face020.jpg"; D="stoßen"; E="face" },
    @{ Row=20; B=123; C="flower/flower024.jpg"; D="spenden"; E="flower" },
    @{ Row=21; B=88; C="face/face017.jpg"; D="angeln"; E="face" },
    @{ Row=22; B=6; C="flower/flower026.jpg"; D="reisen"; E="flower" },
    @{ Row=23; B=24; C="flower/flower010.jpg"; D="kriegen"; E="flower" },
    @{ Row=24; B=115; C="flower/flower023.jpg"; D="planen"; E="flower" },
    @{ Row=25; B=78; C="face/face011.jpg"; D="ärgern"; E="face" },
    @{ Row=26; B=37; C="flower/flower002.jpg"; D="narren"; E="flower" },
    @{ Row=27; B=118; C="flower/flower009.jpg"; D="lügen"; E="flower" },
    @{ Row=28; B=30; C="face/face015.jpg"; D="wecken"; E="face" },
    @{ Row=29; B=25; C="face/face010.jpg"; D="lernen"; E="face" },
    @{ Row=30; B=3; C="face/face016.jpg"; D="nullen"; E="face" },
    @{ Row=31; B=47; C="face/face028.jpg"; D="frischen"; E="face" },
    @{ Row=32; B=106; C="flower/flower007.jpg"; D="bauen"; E="flower" },
    @{ Row=33; B=27; C="flower/flower005.jpg"; D="quellen"; E="flower" }
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 2).Value = $item.B   # B: number
    $ws.Cells.Item($item.Row, 3).Value = $item.C   # C: image
    $ws.Cells.Item($item.Row, 4).Value = $item.D   # D: word
    $ws.Cells.Item($item.Row, 5).Value = $item.E   # E: category
}
